$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new K column values (grades) for the listed rows
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 4
$ws.Range("K10").Value = 5
$ws.Range("K13").Value = 5
$ws.Range("K15").Value = 5
$ws.Range("K17").Value = 4
$ws.Range("K18").Value = 5

# Update the frozen-pane view & active selection to match the saved view state
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("K18").Select()
